# Updates the "cryptos" price/volume table (Sheet1, columns B-E, rows 2-51)
# with refreshed values, matching the automated "cryptos list" GitHub Action
# commit. Columns D (Price) and E (Volume(1h)) hold text-formatted numbers
# (not real numbers), so a handful of cells whose new text would otherwise be
# auto-parsed by Excel as a numeric value (e.g. "216.53", "1.000") are first
# switched to the Text number format ("@") to force them to stay as literal
# strings, exactly like the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Addr='D2'; Val='26.482.70'; Force=$false},
    @{Addr='E2'; Val='  +2.47%  '; Force=$false},
    @{Addr='D3'; Val='1.677.60'; Force=$false},
    @{Addr='E3'; Val='  +3.49%  '; Force=$false},
    @{Addr='E4'; Val='  +0.06%  '; Force=$false},
    @{Addr='D5'; Val='216.53'; Force=$true},
    @{Addr='E5'; Val='  +3.57%  '; Force=$false},
    @{Addr='D6'; Val='0.5319'; Force=$true},
    @{Addr='E6'; Val='  +2.24%  '; Force=$false},
    @{Addr='E7'; Val='  +0.04%  '; Force=$false},
    @{Addr='D8'; Val='0.2675'; Force=$true},
    @{Addr='E8'; Val='  +4.30%  '; Force=$false},
    @{Addr='D9'; Val='0.06392'; Force=$true},
    @{Addr='E9'; Val='  +1.74%  '; Force=$false},
    @{Addr='D10'; Val='21.61'; Force=$true},
    @{Addr='E10'; Val='  +5.92%  '; Force=$false},
    @{Addr='D11'; Val='0.07797'; Force=$true},
    @{Addr='E11'; Val='  +3.64%  '; Force=$false},
    @{Addr='D12'; Val='1.683.90'; Force=$false},
    @{Addr='D13'; Val='4.497'; Force=$true},
    @{Addr='E13'; Val='  +2.92%  '; Force=$false},
    @{Addr='D14'; Val='0.5566'; Force=$true},
    @{Addr='E14'; Val='  +2.05%  '; Force=$false},
    @{Addr='D15'; Val='0.0₅8353'; Force=$false},
    @{Addr='E15'; Val='  +4.79%  '; Force=$false},
    @{Addr='D16'; Val='65.72'; Force=$true},
    @{Addr='E16'; Val='  +2.48%  '; Force=$false},
    @{Addr='D17'; Val='26.523.76'; Force=$false},
    @{Addr='E17'; Val='  +2.64%  '; Force=$false},
    @{Addr='E18'; Val='  +0.02%  '; Force=$false},
    @{Addr='D19'; Val='4.764'; Force=$true},
    @{Addr='E19'; Val='  +2.30%  '; Force=$false},
    @{Addr='D20'; Val='194.69'; Force=$true},
    @{Addr='E20'; Val='  +5.94%  '; Force=$false},
    @{Addr='E21'; Val='  +3.14%  '; Force=$false},
    @{Addr='D22'; Val='6.325'; Force=$true},
    @{Addr='E22'; Val='  +4.26%  '; Force=$false},
    @{Addr='E23'; Val='  +0.06%  '; Force=$false},
    @{Addr='D24'; Val='143.73'; Force=$true},
    @{Addr='E24'; Val='  -0.67%  '; Force=$false},
    @{Addr='D25'; Val='0.1280'; Force=$true},
    @{Addr='E25'; Val='  +6.15%  '; Force=$false},
    @{Addr='D26'; Val='7.433'; Force=$true},
    @{Addr='E26'; Val='  +1.03%  '; Force=$false},
    @{Addr='D27'; Val='16.31'; Force=$true},
    @{Addr='E27'; Val='  +4.87%  '; Force=$false},
    @{Addr='D28'; Val='1.428'; Force=$true},
    @{Addr='E28'; Val='  +5.17%  '; Force=$false},
    @{Addr='E29'; Val='  +4.59%  '; Force=$false},
    @{Addr='D30'; Val='1.275'; Force=$true},
    @{Addr='E30'; Val='  +2.92%  '; Force=$false},
    @{Addr='D31'; Val='3.625'; Force=$true},
    @{Addr='E31'; Val='  +7.14%  '; Force=$false},
    @{Addr='D32'; Val='3.450'; Force=$true},
    @{Addr='E32'; Val='  +3.20%  '; Force=$false},
    @{Addr='D33'; Val='1.690'; Force=$true},
    @{Addr='E33'; Val='  +4.91%  '; Force=$false},
    @{Addr='D34'; Val='1.006'; Force=$true},
    @{Addr='E34'; Val='  +3.59%  '; Force=$false},
    @{Addr='D35'; Val='2.426'; Force=$true},
    @{Addr='D36'; Val='2.784'; Force=$true},
    @{Addr='E36'; Val='  +2.50%  '; Force=$false},
    @{Addr='D37'; Val='0.5738'; Force=$true},
    @{Addr='E37'; Val='  -0.29%  '; Force=$false},
    @{Addr='E38'; Val='  +3.16%  '; Force=$false},
    @{Addr='D39'; Val='6.038'; Force=$true},
    @{Addr='E39'; Val='  +6.71%  '; Force=$false},
    @{Addr='D40'; Val='1.073.12'; Force=$false},
    @{Addr='E40'; Val='  +4.78%  '; Force=$false},
    @{Addr='E41'; Val='  +1.99%  '; Force=$false},
    @{Addr='D42'; Val='1.000'; Force=$true},
    @{Addr='E42'; Val='  -0.21%  '; Force=$false},
    @{Addr='D43'; Val='100.08'; Force=$true},
    @{Addr='E43'; Val='  +0.68%  '; Force=$false},
    @{Addr='D44'; Val='1.825.77'; Force=$false},
    @{Addr='E44'; Val='  +3.21%  '; Force=$false},
    @{Addr='D45'; Val='57.03'; Force=$true},
    @{Addr='E45'; Val='  +4.73%  '; Force=$false},
    @{Addr='B46'; Val='BabyDogeCoin'; Force=$false},
    @{Addr='C46'; Val='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; Force=$false},
    @{Addr='D46'; Val='0.0₈104'; Force=$false},
    @{Addr='E46'; Val='  -4.84%  '; Force=$false},
    @{Addr='B47'; Val='EnergySwap'; Force=$false},
    @{Addr='C47'; Val='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Force=$false},
    @{Addr='D47'; Val='8.145'; Force=$true},
    @{Addr='E47'; Val='  +2.81%  '; Force=$false},
    @{Addr='D48'; Val='1.004'; Force=$true},
    @{Addr='E48'; Val='  +0.39%  '; Force=$false},
    @{Addr='E49'; Val='  +1.00%  '; Force=$false},
    @{Addr='D50'; Val='1.470'; Force=$true},
    @{Addr='E50'; Val='  +6.70%  '; Force=$false},
    @{Addr='D51'; Val='6.031'; Force=$true},
    @{Addr='E51'; Val='  +3.75%  '; Force=$false}
)

foreach ($c in $changes) {
    if ($c.Force) {
        $ws.Range($c.Addr).NumberFormat = "@"
    }
    $ws.Range($c.Addr).Value = $c.Val
}

Write-Host "Applied $($changes.Count) cell changes"